$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=Name, B=PlayerID, C=Last.Updated, D=Injury.Type, E=Injury.Details
# Replace row 2 data with the Jairo Diaz record and remove row 3 (Hanigan) entirely.

$ws.Range("A2").Value = "Jairo Diaz"
$ws.Range("B2").Value = "diazja01"
$ws.Range("C2").Value = "September 10 2017"
$ws.Range("D2").Value = "Undisclosed"
$ws.Range("E2").Value = "Diaz is on the 60-day disabled list with an unknown injury ending his season."

# Delete the old row 3 (Ryan Hanigan) by shifting cells up.
$ws.Range("A3:E3").Delete()

# Row 2 no longer needs the taller custom row height - restore default auto height.
$ws.Rows.Item(2).AutoFit()

$ws.Range("A13").Select()
